# Add a new paragraph ("Threshold for all lists" / line-break / "showcase
# head and tail nodes") right after the paragraph that talks about
# classifying printed students by domestic / international status, and
# right before the "Project Management:" paragraph.

$d = $word.ActiveDocument

$anchor = $d.Content
$found = $anchor.Find.Execute(
    "Another innovation of ours is when all students are printed, it classifies them by domestic / international status ",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if (-not $found) {
    throw "Could not find anchor paragraph text"
}

# Collapse to the end of the found text (end of that paragraph's run) and
# create a brand-new, empty paragraph right after it.
$anchor.Collapse(0) | Out-Null
$anchor.InsertParagraphAfter() | Out-Null

# Move into the newly created paragraph (just past the paragraph mark that
# InsertParagraphAfter left behind) and fill it with the new content. A
# vertical-tab character (Chr 11) is how Word represents a manual line
# break ("Shift+Enter" / <w:br w:type="textWrapping"/>) inside run text.
$anchor.Collapse(0) | Out-Null
$anchor.MoveStart(1, 1) | Out-Null

$lineBreak = [char]11
$newText = "Threshold for all lists" + $lineBreak + "showcase head and tail nodes"
$anchor.InsertAfter($newText) | Out-Null
